$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Project_Admin_Data"

# Write header/data cells in the same order as the original authoring
# (this keeps shared-string indices lined up with the target file)
$ws.Range("A1").Value = "Job code"
$ws.Range("B1").Value = "SOW Fee"
$ws.Range("G1").Value = "Breakfast"
$ws.Range("H1").Value = "Lunch"
$ws.Range("I1").Value = "Dinner"
$ws.Range("J1").Value = "Others"
$ws.Range("A2").Value = "Test1"
$ws.Range("C1").Value = "Speaker Bureau1"
$ws.Range("D1").Value = "Speaker Bureau2"
$ws.Range("E1").Value = "Account Service1"
$ws.Range("F1").Value = "Account Service2"

$ws.Range("B2").Value = 123
$ws.Range("C2").Value = 65
$ws.Range("D2").Value = 65
$ws.Range("E2").Value = 12
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 6

# Bold styling (size 9, Arial, color #444444) applied to C1, D1, H1
$f = $ws.Range("C1").Font
$f.Name = "Arial"
$f.Size = 9
$f.Bold = $true
$f.Color = 4473924

$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column widths (values chosen so the stored/quantized width lands as close
# as possible to the target: 12.5703125, 11.5703125, 17.42578125, 17.42578125,
# 18.28515625, 18.28515625, 12.5703125, 11.28515625, 10.7109375, 11.5703125)
$ws.Columns.Item(1).ColumnWidth = 11.666666666666666
$ws.Columns.Item(2).ColumnWidth = 10.666666666666666
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(4).ColumnWidth = 16.666666666666668
$ws.Columns.Item(5).ColumnWidth = 17.5
$ws.Columns.Item(6).ColumnWidth = 17.5
$ws.Columns.Item(7).ColumnWidth = 11.666666666666666
$ws.Columns.Item(8).ColumnWidth = 10.5
$ws.Columns.Item(9).ColumnWidth = 9.833333333333334
$ws.Columns.Item(10).ColumnWidth = 10.666666666666666

# Page setup (portrait)
$ws.PageSetup.Orientation = 1

# Selection on F2 like in the diff
$ws.Range("F2").Select()
